$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '76.527.99'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '2.906.53'
$ws.Range("E3").Value = '  +7.76%  '
$ws.Range("E4").Value = '  +0.14%  '
Set-TextValue $ws.Cells.Item(5, 4) '197.54'
$ws.Range("E5").Value = '  +5.29%  '
Set-TextValue $ws.Cells.Item(6, 4) '602.31'
$ws.Range("E6").Value = '  +2.23%  '
Set-TextValue $ws.Cells.Item(8, 4) '0.557'
$ws.Range("E8").Value = '  +3.06%  '
$ws.Range("E9").Value = '  -1.28%  '
$ws.Range("D10").Value = '2.905.89'
$ws.Range("E10").Value = '  +7.81%  '
$ws.Range("E11").Value = '  +11.06%  '
Set-TextValue $ws.Cells.Item(12, 4) '0.161'
$ws.Range("E12").Value = '  -1.26%  '
Set-TextValue $ws.Cells.Item(13, 4) '4.96'
$ws.Range("E13").Value = '  +4.91%  '
$ws.Range("D14").Value = '3.416.67'
$ws.Range("E14").Value = '  +7.10%  '
$ws.Range("D15").Value = '76.413.83'
$ws.Range("E15").Value = '  +0.86%  '
Set-TextValue $ws.Cells.Item(16, 4) '27.71'
$ws.Range("E16").Value = '  +4.33%  '
Set-TextValue $ws.Cells.Item(17, 4) '0.0000191'
$ws.Range("D18").Value = '2.909.72'
$ws.Range("E18").Value = '  +7.69%  '
Set-TextValue $ws.Cells.Item(19, 4) '9.05'
$ws.Range("E19").Value = '  -3.23%  '
Set-TextValue $ws.Cells.Item(20, 4) '12.66'
$ws.Range("E20").Value = '  +5.56%  '
Set-TextValue $ws.Cells.Item(21, 4) '385.52'
$ws.Range("E21").Value = '  +2.85%  '
$ws.Range("E22").Value = '  +1.39%  '
Set-TextValue $ws.Cells.Item(23, 4) '4.16'
$ws.Range("E23").Value = '  +1.96%  '
Set-TextValue $ws.Cells.Item(24, 4) '71.94'
$ws.Range("E24").Value = '  +2.68%  '
Set-TextValue $ws.Cells.Item(25, 4) '0.999'
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Cells.Item(26, 4) '4.28'
$ws.Range("E26").Value = '  +2.23%  '
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '3.032.67'
$ws.Range("E27").Value = '  +7.33%  '
Set-TextValue $ws.Cells.Item(28, 4) '9.83'
$ws.Range("E28").Value = '  +4.79%  '
Set-TextValue $ws.Cells.Item(29, 4) '0.0000109'
$ws.Range("E29").Value = '  +14.39%  '
$ws.Range("E30").Value = '  -0.05%  '
Set-TextValue $ws.Cells.Item(31, 4) '1.43'
$ws.Range("E31").Value = '  +0.79%  '
Set-TextValue $ws.Cells.Item(32, 4) '517.97'
$ws.Range("E32").Value = '  -0.90%  '
Set-TextValue $ws.Cells.Item(33, 4) '7.84'
$ws.Range("E33").Value = '  +0.67%  '
$ws.Range("E34").Value = '  +3.93%  '
Set-TextValue $ws.Cells.Item(35, 4) '1.00'
Set-TextValue $ws.Cells.Item(36, 4) '166.68'
$ws.Range("E36").Value = '  +2.39%  '
Set-TextValue $ws.Cells.Item(37, 4) '20.28'
$ws.Range("E37").Value = '  +5.09%  '
Set-TextValue $ws.Cells.Item(38, 4) '0.118'
$ws.Range("E38").Value = '  -1.85%  '
Set-TextValue $ws.Cells.Item(39, 4) '19.73'
$ws.Range("E39").Value = '  +1.68%  '
Set-TextValue $ws.Cells.Item(40, 4) '184.10'
$ws.Range("E40").Value = '  +8.27%  '
Set-TextValue $ws.Cells.Item(42, 4) '0.350'
$ws.Range("E42").Value = '  +5.70%  '
Set-TextValue $ws.Cells.Item(43, 4) '5.13'
$ws.Range("E43").Value = '  +1.82%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Cells.Item(44, 4) '1.70'
$ws.Range("E44").Value = '  -0.74%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Cells.Item(45, 4) '0.0936'
$ws.Range("E45").Value = '  +10.68%  '
$ws.Range("E46").Value = '  +3.70%  '
$ws.Range("E47").Value = '  +2.19%  '
Set-TextValue $ws.Cells.Item(48, 4) '2.39'
$ws.Range("E48").Value = '  +0.63%  '
Set-TextValue $ws.Cells.Item(49, 4) '0.700'
$ws.Range("E49").Value = '  +17.98%  '
Set-TextValue $ws.Cells.Item(50, 4) '0.586'
$ws.Range("E50").Value = '  +8.68%  '
Set-TextValue $ws.Cells.Item(51, 4) '3.79'
$ws.Range("E51").Value = '  +3.61%  '
